# LoadFlow supernode BuildMatrix bug fix
#
# The "Tanh PQ PV" helper sheet models a reactive-power (Q) limiter using a
# TANH curve between Qmin (B2) and Qmax (B3). The supernode BuildMatrix bug
# fix narrows the Q band from [-300, 300] to [-20, 200] and introduces an
# explicit "offs" (B7) term, computed from the new asymmetric band so the
# curve's inflection point lines up correctly:
#   offs = ATANH(Qm / Qs) / Beta
#
# Qs (B5) and Qm (B6) are already formula-driven off B2/B3, so they (and the
# whole E1:E21 / chart-cache dependent values) recompute automatically once
# B2/B3 change and B7 becomes a formula instead of a literal 0.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Qmin: -300 -> -20
$ws.Range("B2").Value = -20

# Qmax: 300 -> 200
$ws.Range("B3").Value = 200

# offs: was a literal 0, now derived from Qm/Qs/Beta
$ws.Range("B7").Formula = "=ATANH(B6/B5)/B4"
